$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "...is opened " | "to the street and features a paved" + " driveway
# and walkway up to the front door."
#   -> "...is opened to the street and features a " | "paved" + " driveway
# and " + "cobble stone " + "walkway up to the front door."
# The paragraph break moves later (after "features a "), and "cobble stone "
# is inserted before "walkway".
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Replacement.ClearFormatting()
$find1 = "is opened ^pto the street and features a paved driveway and walkway up to the front door."
$repl1 = "is opened to the street and features a " + "`r" + "paved driveway and cobble stone walkway up to the front door."
$ok1 = $r1.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)
Write-Output "Edit1: $ok1"

# ---------------------------------------------------------------------------
# Edit 2: Re-flow the paragraph breaks around the Beckett's-front-yard text,
# adding a space where the old paragraph break used to be, and inserting one
# new blank paragraph before "particularly remarkable...".
# ---------------------------------------------------------------------------
$apos = [char]0x2019
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Replacement.ClearFormatting()
$find2 = "The neighborhood children were all drawn to^pthe Beckett" + $apos + "s front yard where they would all play till dusk. Then they would move to the back yard for a couple more hours until it was time to go to their own respective homes. Nobody new why this property ^pwas so special with the kids, there was nothing special or "
$repl2 = "The neighborhood children were all drawn to the Beckett" + $apos + "s front yard where they would all play till dusk. " + "`r" + "Then they would move to the back yard for a couple more hours until it was time to go to their own respective homes. Nobody new why this property was so special with the kids, there was nothing special or "
$ok2 = $r2.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)
Write-Output "Edit2: $ok2"

# Insert a new blank paragraph right before "particularly remarkable...".
$r2b = $d.Content
$found2b = $r2b.Find.Execute("particularly remarkable about the house or the yard", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit2 locate: $found2b"
if ($found2b) {
    $pos2b = $r2b.Start
    $insA = $d.Range($pos2b, $pos2b)
    $insA.InsertParagraphAfter()
    $insB = $d.Range($pos2b, $pos2b)
    $insB.InsertParagraphAfter()
}

# ---------------------------------------------------------------------------
# Edit 3: Merge "...just as " / "long as ..." into one paragraph: "...just as
# long as ...".
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Replacement.ClearFormatting()
$find3 = "Kids are drawn to the unknown just as ^plong as"
$repl3 = "Kids are drawn to the unknown just as long as"
$ok3 = $r3.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)
Write-Output "Edit3: $ok3"

# ---------------------------------------------------------------------------
# Edit 4: Split the trailing " One of the games the children played " into
# " " (stays in the previous paragraph) + a new blank paragraph + "One of the
# games the children played " (new final paragraph).
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("extra air of security. One of the games the children played", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Edit4 locate: $found4"
if ($found4) {
    $marker4 = "extra air of security. "
    $pos4 = $r4.Start + $marker4.Length
    $insC = $d.Range($pos4, $pos4)
    $insC.InsertParagraphAfter()
    $insD = $d.Range($pos4, $pos4)
    $insD.InsertParagraphAfter()
}
